# Update build version/timestamp strings throughout the workbook.
# Old build timestamp: January 30 2026 16.19.47 EST
# New build timestamp: February 02 2026 12.49.33 EST

$oldVersion = "mines - January 30 (built on January 30 2026 16.19.47 EST)"
$newVersion = "mines - January 30 (built on February 02 2026 12.49.33 EST)"

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

# "About" sheet: A2 holds the version banner
$wsAbout.Range("A2").Value = "Version: $newVersion"

# "About" sheet: A6 holds the recommended citation text referencing the version
$wsAbout.Range("A6").Value = "Recommended Citation:  " + [char]34 + "Global Energy Monitor, Coal mine boundaries and methane sources for Foxleigh Coal Mine, Australia, M0040, version '$newVersion'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# "Boundaries and methane sources" sheet: column S (build_version) rows 2-8
for ($row = 2; $row -le 8; $row++) {
    $wsData.Range("S$row").Value = $newVersion
}
